$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the row-8..row-15 "style template" (bordered, bold, centered index
# column) down onto the two brand-new rows (16 and 17) before we touch any
# values, so the new A16/A17 cells pick up the same cell style (s="1") as
# every other row in the index column.
$ws.Cells.Item(15, 1).Copy()
$ws.Cells.Item(16, 1).PasteSpecial(-4122)
$ws.Cells.Item(16, 1).Copy()
$ws.Cells.Item(17, 1).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Two new contingency entries ("line7"/"line8") were inserted into the
# table right after "line6", pushing all of the "extr#" rows down by two
# rows. Rewrite rows 8-17 (data rows 7-16) of the table with their final
# values.

# row 8: line7
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

# row 9: line8
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# row 10: extr1
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "extr1"
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

# row 11: extr2
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "extr2"
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = $true

# row 12: extr3
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "extr3"
$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 4).Value = 11
$ws.Cells.Item(12, 5).Value = $true

# row 13: extr4
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "extr4"
$ws.Cells.Item(13, 3).Value = 7
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = $false

# row 14: extr5
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "extr5"
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = $false

# row 15: extr6
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "extr6"
$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $true

# row 16: extr7 (new row)
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "extr7"
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $true

# row 17: extr8 (new row)
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "extr8"
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $true
